$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.471.53'
$ws.Range("E2").Value = '  -1.93%  '
$ws.Range("D3").Value = '2.629.60'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.84'
$ws.Range("E5").Value = '  -3.43%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.07'
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.646'
$ws.Range("E7").Value = '  +4.90%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -5.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.79'
$ws.Range("E10").Value = '  -0.36%  '
$ws.Range("E11").Value = '  -2.51%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.50'
$ws.Range("E13").Value = '  -1.86%  '
$ws.Range("E14").Value = '  -6.67%  '
$ws.Range("D15").Value = '3.103.11'
$ws.Range("E15").Value = '  -1.38%  '
$ws.Range("D16").Value = '64.282.44'
$ws.Range("E16").Value = '  -1.99%  '
$ws.Range("D17").Value = '2.631.74'
$ws.Range("E17").Value = '  -1.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.24'
$ws.Range("E18").Value = '  -3.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.68'
$ws.Range("E19").Value = '  -2.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.38'
$ws.Range("E20").Value = '  -1.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '346.02'
$ws.Range("E21").Value = '  -1.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  -2.04%  '
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("E25").Value = '  +2.90%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.36'
$ws.Range("E26").Value = '  -3.32%  '
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '557.57'
$ws.Range("E28").Value = '  +4.37%  '
$ws.Range("E29").Value = '  -2.15%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("E32").Value = '  -3.06%  '
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.59'
$ws.Range("E34").Value = '  +2.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.27'
$ws.Range("E35").Value = '  -3.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.411'
$ws.Range("E36").Value = '  -2.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.02'
$ws.Range("E37").Value = '  -2.57%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  -0.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '154.63'
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.44'
$ws.Range("E42").Value = '  +4.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '158.64'
$ws.Range("E43").Value = '  -2.38%  '
$ws.Range("E44").Value = '  -3.06%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0597'
$ws.Range("E45").Value = '  -2.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.75'
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.635'
$ws.Range("E48").Value = '  +3.14%  '
$ws.Range("E49").Value = '  -2.81%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.08'
$ws.Range("E50").Value = '  -4.12%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0239'
$ws.Range("E51").Value = '  -6.21%  '
